$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 362 (shifts existing row 362 and below down by one)
$ws.Rows("362:362").Insert()

# Populate the new row's cells
$ws.Range("A362").Value = "Mary McLeod Bethune School"
$ws.Range("B362").Value = "Mary M. Bethune School"

# Update the view to match target state
$ws.Application.ActiveWindow.ScrollRow = 348
$ws.Range("B363").Select()
$ws.Application.ActiveWindow.RangeSelection.Value | Out-Null
